$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A84").NumberFormat = "@"
$ws.Range("A84").Value = "2025-10-17"
$ws.Range("A84").Style = "Normal"
$ws.Range("B84").Value = "21:22:34"
$ws.Range("C84").Value = "1.00 EUR = 1,709.2239"
